# Sprint 3 Präsentation.pptx - apply authored edits
#
# 1) Slide 1 title: merge the "Ziele der Gruppe nach " + "S" runs into a
#    single run "Ziele der Gruppe nach S" (text content stays identical,
#    only the run split changes).
# 2) Slide 3 "Eingabe/Ausgabe" textbox: split "Gewicht eines Paketes" into
#    "Gewicht " + "pro Packung", and "Preis" into "Preis " + "pro Packung".

$p = $ppt.ActivePresentation

# --- Slide 1: Title run merge -------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleTr = $slide1.Shapes.Item(1).TextFrame.TextRange

$run1 = $titleTr.Runs(1)
$run2 = $titleTr.Runs(2)
$run1.Text = "Ziele der Gruppe nach S"
$run2.Text = ""

# --- Slide 3: "Gewicht eines Paketes" -> "Gewicht " + "pro Packung" -----
$slide3 = $p.Slides.Item(3)
$boxTr = $slide3.Shapes.Item(6).TextFrame.TextRange

$gewichtPara = $boxTr.Paragraphs(8)
$gewichtRun = $gewichtPara.Runs(1)
$gewichtRun.Text = "Gewicht "
$gewichtRun.InsertAfter("pro Packung") | Out-Null

# --- Slide 3: "Preis" -> "Preis " + "pro Packung" ------------------------
$preisPara = $boxTr.Paragraphs(9)
$preisRun = $preisPara.Runs(1)
$preisRun.Text = "Preis "
$preisRun.InsertAfter("pro Packung") | Out-Null
